$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "User" entity table (A2:B8) documents its fields in column A with a
# flags column in B (e.g. "Unique"). Mark the UserName and Password fields
# as unique, matching how the Customer table already flags its unique
# fields (Email/PhoneNumber) in column E.
$ws.Range("B4").Value = "Unique"
$ws.Range("B5").Value = "Unique"

# Update the sheet's view: no longer frozen/scrolled to row 7, and the
# active selection moves to B7 instead of C17.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B7").Select()
